# Generate Report for Handoff
# Swap the "baacd752-..." and "c1e58909-..." file rows on the Overview,
# zh-cn and de-de sheets: baacd752 moves from "Handed back" (row 2) to
# "Ready for handoff" (row 3) with fresh handoff timestamps and an error
# detail message, while c1e58909 takes over the "Handed back: in sync
# with en-US" state in row 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.md"
$ov.Range("B2").Value = "e2e\c1e58909-a01c-40a0-a0c9-266b2875f041.md"

$ov.Range("A3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.md"
$ov.Range("B3").Value = "e2e\baacd752-d650-4f23-89e6-0db079eccae6.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-17 16:47:30"

# Recreate the two hyperlinks, keeping the same link targets (rId2 ->
# baacd752 github blob, rId3 -> c1e58909 github blob) but with the
# display text swapped to match the new row contents.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md", "", "", "e2e\c1e58909-a01c-40a0-a0c9-266b2875f041.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/c1e58909-a01c-40a0-a0c9-266b2875f041.md", "", "", "e2e\baacd752-d650-4f23-89e6-0db079eccae6.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.md"
$zh.Range("G2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.e8562b93bdd49870a7773764ab055171c7c4c662.zh-cn.xlf"
$zh.Range("I2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.md"
$zh.Range("J2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.e8562b93bdd49870a7773764ab055171c7c4c662.zh-cn.xlf"

$zh.Range("A3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.3eac469c0f8383a76f0040ae99bcc80501f82a8c.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-17 16:47:24"
$zh.Range("I3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.md"
$zh.Range("J3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.3eac469c0f8383a76f0040ae99bcc80501f82a8c.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00d2cc5ca7e39defbc49aa5424ebb95086668bf6/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md."

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md", "", "", "c1e58909-a01c-40a0-a0c9-266b2875f041.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f580eaca2571a0eec38cb2e5e8a13377293027cc/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md", "", "", "c1e58909-a01c-40a0-a0c9-266b2875f041.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/c1e58909-a01c-40a0-a0c9-266b2875f041.md", "", "", "baacd752-d650-4f23-89e6-0db079eccae6.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f580eaca2571a0eec38cb2e5e8a13377293027cc/e2e/c1e58909-a01c-40a0-a0c9-266b2875f041.md", "", "", "baacd752-d650-4f23-89e6-0db079eccae6.md")

$zh.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.md"
$de.Range("G2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.e8562b93bdd49870a7773764ab055171c7c4c662.de-de.xlf"
$de.Range("I2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.md"
$de.Range("J2").Value = "c1e58909-a01c-40a0-a0c9-266b2875f041.e8562b93bdd49870a7773764ab055171c7c4c662.de-de.xlf"

$de.Range("A3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.3eac469c0f8383a76f0040ae99bcc80501f82a8c.de-de.xlf"
$de.Range("H3").Value = "2016-08-17 16:47:30"
$de.Range("I3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.md"
$de.Range("J3").Value = "baacd752-d650-4f23-89e6-0db079eccae6.3eac469c0f8383a76f0040ae99bcc80501f82a8c.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/00d2cc5ca7e39defbc49aa5424ebb95086668bf6/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md."

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md", "", "", "c1e58909-a01c-40a0-a0c9-266b2875f041.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1caabd6f73b142830c2d6c64fd7e782497b18d45/e2e/baacd752-d650-4f23-89e6-0db079eccae6.md", "", "", "c1e58909-a01c-40a0-a0c9-266b2875f041.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c88c836f7dda94a04d244af0c335add5df96acbd/e2e/c1e58909-a01c-40a0-a0c9-266b2875f041.md", "", "", "baacd752-d650-4f23-89e6-0db079eccae6.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1caabd6f73b142830c2d6c64fd7e782497b18d45/e2e/c1e58909-a01c-40a0-a0c9-266b2875f041.md", "", "", "baacd752-d650-4f23-89e6-0db079eccae6.md")

$de.Columns.Item(16).ColumnWidth = 39.17
